# Updates cryptos list values (price/volume, and a few coin row swaps)
# per the Sat Jul 13 07:55:53 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.100.55"
$ws.Range("E2").Value = "  +2.00%  "

# Row 3
$ws.Range("D3").Value = "3.137.21"
$ws.Range("E3").Value = "  +2.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'535.39"
$ws.Range("E5").Value = "  +2.86%  "

# Row 6
$ws.Range("D6").Value = "'139.17"
$ws.Range("E6").Value = "  +2.99%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  +12.66%  "

# Row 9
$ws.Range("D9").Value = "'7.31"
$ws.Range("E9").Value = "  -0.11%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.109"
$ws.Range("E10").Value = "  +3.09%  "

# Row 11
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.424"
$ws.Range("E11").Value = "  +6.88%  "

# Row 12
$ws.Range("E12").Value = "  +3.22%  "

# Row 13
$ws.Range("D13").Value = "3.677.86"
$ws.Range("E13").Value = "  +2.25%  "

# Row 14
$ws.Range("D14").Value = "'25.77"
$ws.Range("E14").Value = "  +2.50%  "

# Row 15
$ws.Range("E15").Value = "  +5.35%  "

# Row 16
$ws.Range("D16").Value = "58.213.46"
$ws.Range("E16").Value = "  +2.10%  "

# Row 17
$ws.Range("D17").Value = "'6.26"
$ws.Range("E17").Value = "  +6.96%  "

# Row 18
$ws.Range("D18").Value = "3.145.08"
$ws.Range("E18").Value = "  +2.63%  "

# Row 19
$ws.Range("E19").Value = "  +4.77%  "

# Row 20
$ws.Range("D20").Value = "'8.20"
$ws.Range("E20").Value = "  +5.11%  "

# Row 21
$ws.Range("D21").Value = "'376.63"
$ws.Range("E21").Value = "  +8.13%  "

# Row 22
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.07%  "

# Row 23
$ws.Range("D23").Value = "'5.74"
$ws.Range("E23").Value = "  -0.76%  "

# Row 24
$ws.Range("D24").Value = "'70.47"
$ws.Range("E24").Value = "  +2.68%  "

# Row 25
$ws.Range("D25").Value = "'0.516"
$ws.Range("E25").Value = "  +3.90%  "

# Row 26
$ws.Range("E26").Value = "  +0.76%  "

# Row 27
$ws.Range("E27").Value = "  +0.06%  "

# Row 28
$ws.Range("D28").Value = "'8.07"
$ws.Range("E28").Value = "  +12.43%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0882"
$ws.Range("E29").Value = "  +2.03%  "

# Row 30
$ws.Range("D30").Value = "'6.20"
$ws.Range("E30").Value = "  +6.68%  "

# Row 31
$ws.Range("E31").Value = "  +1.95%  "

# Row 32
$ws.Range("D32").Value = "'21.75"
$ws.Range("E32").Value = "  +4.33%  "

# Row 33
$ws.Range("E33").Value = "  +7.37%  "

# Row 34
$ws.Range("E34").Value = "  +4.10%  "

# Row 35
$ws.Range("D35").Value = "'161.79"
$ws.Range("E35").Value = "  +1.61%  "

# Row 36
$ws.Range("E36").Value = "  +5.29%  "

# Row 37
$ws.Range("D37").Value = "'1.35"
$ws.Range("E37").Value = "  +10.60%  "

# Row 38
$ws.Range("D38").Value = "'25.62"
$ws.Range("E38").Value = "  +1.29%  "

# Row 39
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.639.81"
$ws.Range("E39").Value = "  +10.51%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.67"
$ws.Range("E40").Value = "  +5.95%  "

# Row 41
$ws.Range("E41").Value = "  +5.77%  "

# Row 42
$ws.Range("E42").Value = "  +3.08%  "

# Row 43
$ws.Range("E43").Value = "  +6.49%  "

# Row 44
$ws.Range("D44").Value = "'0.700"
$ws.Range("E44").Value = "  +1.58%  "

# Row 45
$ws.Range("D45").Value = "'0.0274"
$ws.Range("E45").Value = "  +5.12%  "

# Row 46
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("E47").Value = "  +4.91%  "

# Row 48
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'0.976"
$ws.Range("E48").Value = "  +2.86%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  +11.25%  "

# Row 50
$ws.Range("D50").Value = "'20.31"
$ws.Range("E50").Value = "  +3.78%  "

# Row 51
$ws.Range("D51").Value = "'0.747"
$ws.Range("E51").Value = "  -0.41%  "

